# TestNG package created 26-July
$wb = $excel.ActiveWorkbook

# Rename "Sheet3" to "url"
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Name = "url"

# Fill column A of the "url" sheet with a header and 15 numbers
$ws3.Range("A1").Value = "url"
for ($i = 1; $i -le 15; $i++) {
    $ws3.Cells.Item($i + 1, 1).Value = $i
}

# Apply the bordered style (style index 1 in the original workbook) to A1:A16
# by copying the format from an existing cell that already carries that style.
$ws1 = $wb.Worksheets.Item("Home")
$ws1.Range("B1").Copy()
$ws3.Range("A1:A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Select A1:A16 on the "url" sheet and make it the active sheet/tab
$ws3.Range("A1:A16").Select()
$ws3.Activate()
